$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the column-dictionary labels (column B) that were renamed.
# The assignment order below matches the order the author's edit introduced
# the new labels, so the shared-string table is rebuilt with the same layout.
$ws.Range("B46").Value = "prior_births_dead"
$ws.Range("B47").Value = "prior_births_living"
$ws.Range("B48").Value = "prior_terminations"
$ws.Range("B49").Value = "prepreg_weight"
$ws.Range("B50").Value = "delivery_method"
$ws.Range("B54").Value = "infant_sex"
$ws.Range("B42").Value = "payment"
$ws.Range("B39").Value = "infections"
$ws.Range("B33").Value = "m_race15"
$ws.Range("B34").Value = "m_race31"
$ws.Range("B36").Value = "m_race6"
$ws.Range("B17").Value = "f_race15"
$ws.Range("B18").Value = "f_race31"
$ws.Range("B19").Value = "f_race6"
$ws.Range("B31").Value = "m_hispanic"
$ws.Range("B16").Value = "f_hispanic"
$ws.Range("B40").Value = "m_morbidity"
$ws.Range("B41").Value = "riskf"

# The shortened text in rows 47 and 49 no longer wraps to two lines, so
# Excel recalculates (shrinks) those row heights back to the sheet default.
$ws.Rows.Item(47).AutoFit()
$ws.Rows.Item(49).AutoFit()

# Restore the view/selection state recorded for the sheet.
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 28
$null = $ws.Range("B41").Select()
